# Applies the update described by the commit "Update gh-pages to output
# generated at 456a3b4" to the 杭州-漫展信息 workbook.
#
# Sheets (in workbook order):
#   1 = 展览      (Exhibitions)
#   2 = 演出      (Performances)
#   3 = 本地生活   (Local life)
#   4 = 全部类型   (All categories, combined view)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Simple "want-to-go" counter refreshes (F column) across existing rows.
$ws1.Range("F3").Value  = 2567
$ws1.Range("F7").Value  = 1364
$ws1.Range("F8").Value  = 1721
$ws1.Range("F9").Value  = 190
$ws1.Range("F11").Value = 2394
$ws1.Range("F12").Value = 505
$ws1.Range("F13").Value = 178
$ws1.Range("F14").Value = 51
$ws1.Range("F17").Value = 97
$ws1.Range("F18").Value = 8717
$ws1.Range("F20").Value = 6802
$ws1.Range("F21").Value = 10986
$ws1.Range("F23").Value = 188
$ws1.Range("F24").Value = 220
$ws1.Range("F25").Value = 302
$ws1.Range("F27").Value = 2419
$ws1.Range("F28").Value = 205
$ws1.Range("F29").Value = 181
$ws1.Range("F30").Value = 2252
$ws1.Range("F31").Value = 327
$ws1.Range("F32").Value = 25
$ws1.Range("F33").Value = 4457
$ws1.Range("F34").Value = 524

# Venue text amended (501 city square addition).
$ws1.Range("D7").Value = "景昙路9号 杭州大厦501城市广场"

# Insert a brand-new row so that the "谢莹内场票" listing (previously row
# 35) moves down to row 36, the "华盟次元" listing (previously row 36)
# moves down to row 37, and row 35 is freed up for a new "陈珂内场票"
# listing.
$ws1.Rows.Item(36).Insert()

# The insert copies row 35's formatting down, but column A needs the
# bordered/centered numbering style - pull it back from row 35.
$ws1.Range("A35").Copy($ws1.Range("A36"))

# New row 36 = the former row 35 contents ("谢莹内场票"), with its
# updated want-to-go count.
$ws1.Range("A36").Value = 35
$ws1.Range("B36").NumberFormat = "@"
$ws1.Range("B36").Value = "2024-12-15"
$ws1.Range("C36").Value = "杭州·AD05动漫展.谢莹内场票"
$ws1.Range("D36").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws1.Range("E36").Value = "2024.12.15 09:00-12.15 17:00"
$ws1.Range("F36").Value = 12
$ws1.Range("G36").Value = 108
$ws1.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=93285"
$ws1.Range("I36").Value = "//i1.hdslb.com/bfs/openplatform/202410/55HRoWBV1728461646171.png"

# Row 35 now becomes the new "陈珂内场票" listing.
$ws1.Range("B35").NumberFormat = "@"
$ws1.Range("B35").Value = "2024-12-14"
$ws1.Range("C35").Value = "杭州·AD05动漫展.陈珂内场票"
$ws1.Range("E35").Value = "2024.12.14 09:00-12.14 17:00"
$ws1.Range("F35").Value = 15
$ws1.Range("G35").Value = 188
$ws1.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=93306"
$ws1.Range("I35").Value = "//i1.hdslb.com/bfs/openplatform/202410/9mZW4Sep1728620496423.jpeg"

# Row 37 is the shifted-down "华盟次元" listing - its running index and
# want-to-go count both changed (index 35 -> 36, count 447 -> 452).
$ws1.Range("A37").Value = 36
$ws1.Range("F37").Value = 452

# ---------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value = 1178

# ---------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F5").Value = 70

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (combined view - mirrors the other sheets' updates)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F5").Value  = 70
$ws4.Range("F7").Value  = 2567
$ws4.Range("D11").Value = "景昙路9号 杭州大厦501城市广场"
$ws4.Range("F11").Value = 1364
$ws4.Range("F13").Value = 1721
$ws4.Range("F15").Value = 190
$ws4.Range("F17").Value = 505
$ws4.Range("F18").Value = 178
$ws4.Range("F19").Value = 51
$ws4.Range("F22").Value = 97
$ws4.Range("F23").Value = 8717
$ws4.Range("F25").Value = 6802
$ws4.Range("F26").Value = 10986
$ws4.Range("F29").Value = 188
$ws4.Range("F30").Value = 220
$ws4.Range("F31").Value = 302
$ws4.Range("F36").Value = 205
$ws4.Range("F37").Value = 25
$ws4.Range("F38").Value = 4457
$ws4.Range("F45").Value = 452

$wb.Save()
